$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-10 Monday" "2025-02-11 Tuesday"

Replace-Text "748÷5=149, 3" "321÷3=107, 0"
Replace-Text "315÷4=78, 3" "135÷2=67, 1"
Replace-Text "339÷8=42, 3" "182÷7=26, 0"
Replace-Text "812÷9=90, 2" "130÷7=18, 4"
Replace-Text "356÷4=89, 0" "150÷4=37, 2"

Replace-Text "344÷3=114, 2" "916÷7=130, 6"
Replace-Text "862÷9=95, 7" "622÷7=88, 6"
Replace-Text "263÷2=131, 1" "158÷5=31, 3"
Replace-Text "857÷8=107, 1" "219÷8=27, 3"
Replace-Text "410÷9=45, 5" "359÷6=59, 5"

Replace-Text "744÷2=372, 0" "434÷4=108, 2"
Replace-Text "588÷8=73, 4" "140÷9=15, 5"
Replace-Text "949÷5=189, 4" "578÷3=192, 2"
Replace-Text "945÷4=236, 1" "150÷6=25, 0"
Replace-Text "695÷7=99, 2" "819÷3=273, 0"

Replace-Text "400÷3=133, 1" "288÷2=144, 0"
Replace-Text "721÷9=80, 1" "328÷5=65, 3"
Replace-Text "184÷5=36, 4" "205÷2=102, 1"
Replace-Text "464÷2=232, 0" "239÷9=26, 5"
Replace-Text "430÷4=107, 2" "372÷3=124, 0"

Replace-Text "897÷6=149, 3" "939÷7=134, 1"
Replace-Text "177÷9=19, 6" "278÷9=30, 8"
Replace-Text "597÷7=85, 2" "946÷9=105, 1"
Replace-Text "852÷4=213, 0" "634÷5=126, 4"
Replace-Text "864÷5=172, 4" "347÷6=57, 5"

Write-Output "Done"
